# Normalize the "Recorded By" (column G) entries on the Session Analysis
# Results sheet: for every multi-value, comma-separated "Recorded By" list
# that doesn't already start with the human reviewer's email
# (dnasr281@gmail.com), reverse the order of the entries.
#
# Single-value cells, and cells already led by dnasr281@gmail.com, are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }
    if ($value -notlike "*,*") { continue }

    $parts = $value -split ", "
    if ($parts.Count -le 1) { continue }

    $firstPart = $parts[0].Trim()
    if ($firstPart -eq "dnasr281@gmail.com") { continue }

    $reversedParts = $parts[($parts.Count - 1)..0]
    $cell.Value2 = [string]::Join(", ", $reversedParts)
}
